{"js": "// Replace \"AR ${AR}\" (heading) with \"Recommendation ${REC}\",\n// and \" for this AR will be\" with \" for this recommendation will be\"\n// (commit: \"Getting rid of AR\" / \"All ARs are replaced with recommendations\").\n\nconst body = context.document.body;\n\n// --- Change 1: the document heading \"AR ${AR}: ...\" -> \"Recommendation ${REC}: ...\" ---\nconst headingResults = body.search(\"AR ${AR}\", { matchCase: true, matchWholeWord: false });\nheadingResults.load(\"items\");\nawait context.sync();\n\nfor (const r of headingResults.items) {\n  r.insertText(\"Recommendation ${REC}\", \"Replace\");\n}\nawait context.sync();\n\n// --- Change 2: \" for this AR will be\" -> \" for this recommendation will be\" ---\nconst bodyResults = body.search(\" for this AR will be\", { matchCase: true, matchWholeWord: false });\nbodyResults.load(\"items\");\nawait context.sync();\n\nfor (const r of bodyResults.items) {\n  r.insertText(\" for this recommendation will be\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Replace \"AR ${AR}\" (heading) with \"Recommendation ${REC}\",\n# and \" for this AR will be\" with \" for this recommendation will be\"\n# (commit: \"Getting rid of AR\" / \"All ARs are replaced with recommendations\").\n\n$d = $word.ActiveDocument\n\n# --- Change 1: the document heading \"AR ${AR}: ...\" -> \"Recommendation ${REC}: ...\" ---\n$find1 = $d.Content.Find\n$find1.Text = 'AR ${AR}'\n$find1.Replacement.Text = 'Recommendation ${REC}'\n$find1.Execute([ref]'AR ${AR}', $false, $false, $false, $false, $false, $true, 1, $false, [ref]'Recommendation ${REC}', 2)\n\n# --- Change 2: \" for this AR will be\" -> \" for this recommendation will be\" ---\n$find2 = $d.Content.Find\n$find2.Text = ' for this AR will be'\n$find2.Replacement.Text = ' for this recommendation will be'\n$find2.Execute([ref]' for this AR will be', $false, $false, $false, $false, $false, $true, 1, $false, [ref]' for this recommendation will be', 2)\n"}
